# Update cryptos list: refresh D (Price) / E (Volume 1h) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.519.60"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.10"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.16"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.17"
$ws.Range("E8").Value = "  +5.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.07"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -1.84%  "
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.794.63"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.561.17"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.479.42"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.15"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.45"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("E21").Value = "  -3.08%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -5.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.11"
$ws.Range("E25").Value = "  +7.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.78"
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.04"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.41"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0480"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("E32").Value = "  -3.75%  "
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.08"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.391.42"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.64"
$ws.Range("E39").Value = "  +4.61%  "
$ws.Range("E40").Value = "  -1.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.521"
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.788"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0468"
$ws.Range("E45").Value = "  +3.63%  "
$ws.Range("E46").Value = "  -5.15%  "
$ws.Range("E47").Value = "  -1.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "62.85"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.707.73"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.34"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0101"
$ws.Range("E51").Value = "  -2.28%  "
